$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.264.88"
$ws.Range("D2").Style = $dStyle
$ws.Range("E2").Value = "  +1.19%  "

$dStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.857.41"
$ws.Range("D3").Style = $dStyle
$ws.Range("E3").Value = "  +1.80%  "

$dStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = $dStyle
$ws.Range("E4").Value = "  -0.34%  "

$dStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.11"
$ws.Range("D5").Style = $dStyle
$ws.Range("E5").Value = "  +1.14%  "

$dStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = $dStyle
$ws.Range("E6").Value = "  -0.21%  "

$dStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4654"
$ws.Range("D7").Style = $dStyle
$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("E8").Value = "  +0.56%  "

$dStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07300"
$ws.Range("D9").Style = $dStyle
$ws.Range("E9").Value = "  -0.44%  "

$dStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8917"
$ws.Range("D10").Style = $dStyle
$ws.Range("E10").Value = "  +1.91%  "

$dStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.08"
$ws.Range("D11").Style = $dStyle
$ws.Range("E11").Value = "  +2.14%  "

$dStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07880"
$ws.Range("D12").Style = $dStyle
$ws.Range("E12").Value = "  +0.26%  "

$dStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.839.43"
$ws.Range("D13").Style = $dStyle
$ws.Range("E13").Value = "  -0.03%  "

$dStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.417"
$ws.Range("D14").Style = $dStyle
$ws.Range("E14").Value = "  +1.59%  "

$dStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.518"
$ws.Range("D15").Style = $dStyle
$ws.Range("E15").Value = "  -0.44%  "

$dStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.64"
$ws.Range("D16").Style = $dStyle
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("E17").Value = "  -0.30%  "

$dStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008938"
$ws.Range("D18").Style = $dStyle
$ws.Range("E18").Value = "  +1.35%  "

$dStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = $dStyle
$ws.Range("E19").Value = "  -0.14%  "

$ws.Range("E20").Value = "  -0.06%  "

$dStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.293.76"
$ws.Range("D21").Style = $dStyle
$ws.Range("E21").Value = "  +1.21%  "

$dStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.081"
$ws.Range("D22").Style = $dStyle
$ws.Range("E22").Value = "  -0.33%  "

$ws.Range("E23").Value = "  +0.12%  "

$dStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.070.18"
$ws.Range("D24").Style = $dStyle
$ws.Range("E24").Value = "  +2.15%  "

$dStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.037"
$ws.Range("D25").Style = $dStyle
$ws.Range("E25").Value = "  +9.64%  "

$dStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.73"
$ws.Range("D26").Style = $dStyle
$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("E27").Value = "  +0.11%  "

$dStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.050"
$ws.Range("D28").Style = $dStyle
$ws.Range("E28").Value = "  +0.76%  "

$dStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.94"
$ws.Range("D29").Style = $dStyle
$ws.Range("E29").Value = "  +0.43%  "

$dStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.048"
$ws.Range("D30").Style = $dStyle
$ws.Range("E30").Value = "  -1.00%  "

$dStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08846"
$ws.Range("D31").Style = $dStyle
$ws.Range("E31").Value = "  -0.29%  "

$dStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.147"
$ws.Range("D32").Style = $dStyle
$ws.Range("E32").Value = "  +6.39%  "

$dStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7710"
$ws.Range("D33").Style = $dStyle
$ws.Range("E33").Value = "  +5.71%  "

$dStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.171"
$ws.Range("D34").Style = $dStyle
$ws.Range("E34").Value = "  +3.60%  "

$dStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.527"
$ws.Range("D35").Style = $dStyle
$ws.Range("E35").Value = "  +2.11%  "

$dStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.703"
$ws.Range("D36").Style = $dStyle
$ws.Range("E36").Value = "  +9.86%  "

$ws.Range("E37").Value = "  +3.05%  "

$dStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01943"
$ws.Range("D38").Style = $dStyle
$ws.Range("E38").Value = "  +0.23%  "

$dStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05230"
$ws.Range("D39").Style = $dStyle
$ws.Range("E39").Value = "  +0.33%  "

$dStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.952"
$ws.Range("D40").Style = $dStyle
$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("E41").Value = "  -0.44%  "

$dStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5122"
$ws.Range("D42").Style = $dStyle
$ws.Range("E42").Value = "  -0.17%  "

$ws.Range("E43").Value = "  +0.40%  "

$dStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.526"
$ws.Range("D44").Style = $dStyle
$ws.Range("E44").Value = "  +4.72%  "

$dStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4800"
$ws.Range("D45").Style = $dStyle
$ws.Range("E45").Value = "  -0.31%  "

$dStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.38"
$ws.Range("D46").Style = $dStyle
$ws.Range("E46").Value = "  +2.36%  "

$dStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").Style = $dStyle
$ws.Range("E47").Value = "  -0.22%  "

$dStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.85"
$ws.Range("D48").Style = $dStyle
$ws.Range("E48").Value = "  +1.10%  "

$dStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.646"
$ws.Range("D49").Style = $dStyle
$ws.Range("E49").Value = "  +1.66%  "

$dStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06204"
$ws.Range("D50").Style = $dStyle
$ws.Range("E50").Value = "  +0.05%  "

$dStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.49"
$ws.Range("D51").Style = $dStyle
$ws.Range("E51").Value = "  +1.52%  "
